$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Cas de test 1"

# Row 2
$ws.Range("E2").Value = "Cas de test 1"
$ws.Range("F2").Value = 3.668316228261497
$ws.Range("I2").Value = 5.356121928268019

# Row 3 (was amazon.fr/, now ausy.fr/fr/) + new E3/F3 for "Cas de test 2"
$ws.Range("A3").Value = "https://www.ausy.fr/fr/"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 3.151041094044596
$ws.Range("E3").Value = "Cas de test 2"
$ws.Range("F3").Value = 1.687805700006522

# Row 4
$ws.Range("A4").Value = "https://www.ausy.fr/fr/carrieres/"
$ws.Range("B4").Value = 71
$ws.Range("C4").Value = 3.415620132068172

# Row 5
$ws.Range("A5").Value = "https://www.ausy.fr/fr/carrieres/toutes-nos-offres/q-auto/"
$ws.Range("B5").Value = 73
$ws.Range("C5").Value = 3.668316228261497

# Row 6 previously held amazon order-history data; it is removed entirely,
# and new rows 7-12 are appended below with the second test case block.
$ws.Range("A6").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# Row 7
$ws.Range("A7").Value = "Cas de test 2"

# Row 8 (headers)
$ws.Range("A8").Value = "URL"
$ws.Range("B8").Value = "Propreté (en %)"
$ws.Range("C8").Value = "Consommation de CO2 (en g)"

# Row 9
$ws.Range("A9").Value = "https://www.amazon.fr/"
$ws.Range("B9").Value = 36
$ws.Range("C9").Value = 0.7976718296779319

# Row 10
$ws.Range("A10").Value = "https://www.amazon.fr/deals?ref_=nav_cs_gb"
$ws.Range("B10").Value = 66
$ws.Range("C10").Value = 1.150309954035468

# Row 11
$ws.Range("A11").Value = "https://www.amazon.fr/deal/3a51f27b?showVariations=true&pf_rd_r=AW6MJV6KCYCK3YT21XYQ&pf_rd_t=Events&pf_rd_i=deals&pf_rd_p=3c3f3ff2-f80e-428b-aff0-b0531c852487&pf_rd_s=slot-14&ref=dlx_deals_gd_dcl_img_1_3a51f27b_dt_sl14_87"
$ws.Range("B11").Value = 71
$ws.Range("C11").Value = 1.452404730597977

# Row 12
$ws.Range("A12").Value = "https://www.amazon.fr/gp/your-account/order-history?ref_=ya_d_c_yo"
$ws.Range("B12").Value = 79
$ws.Range("C12").Value = 1.687805700006522
